# Semana 38 de 2025: add week-38 column (AO) to the weekly IRA report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AO1 holds the week number as text ("38"), matching the
# other week-header cells (D1.."AN1" = "1".."37"). A leading apostrophe
# forces Excel to store it as text instead of a number.
$ws.Range("AO1").Formula = "'38"

# Fix for row 35 (week 37 total revised upward) plus its new week-38 value.
$ws.Range("AN35").Value = 39
$ws.Range("AO35").Value = 37

# New week-38 counts (column AO) for every UPGD row that reported data.
$weekValues = @{
    2  = 74
    5  = 1
    6  = 68
    7  = 29
    8  = 26
    9  = 3
    10 = 2
    11 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    23 = 5
    25 = 46
    28 = 255
    29 = 0
    31 = 2
    36 = 2
    37 = 9
    38 = 86
    41 = 16
    42 = 5
    43 = 33
    44 = 175
    45 = 138
    46 = 149
    47 = 3
    48 = 93
    49 = 2
    50 = 0
    51 = 6
    52 = 5
    53 = 4
    54 = 1
    55 = 0
    56 = 3
    57 = 77
    58 = 6
}

foreach ($row in $weekValues.Keys) {
    $ws.Cells.Item($row, 41).Value = $weekValues[$row]
}
